# Populate the "usuarios" sheet with the owner/admin record created by the
# new "create owner from Excel" feature: Antonio / ctr / admin-Antonio /
# false / crt2 across A1:E1, leaving row 2 as blank placeholder cells
# (matching the pre-existing blank A2/B2 pattern).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Antonio"
$ws.Range("B1").Value = "ctr"
$ws.Range("C1").Value = "admin-Antonio"
$ws.Range("E1").Value = "crt2"

# D1 must hold the literal TEXT "false" (not a boolean). Assigning the
# string directly gets auto-coerced to a logical value, so compute it as a
# formula result and paste back as a value, which preserves its text type.
$ws.Range("G1").Formula = "=""false"""
$ws.Range("G1").Copy()
$ws.Range("D1").PasteSpecial(-4163)
$ws.Range("G1").ClearContents()

# Extend row 2 with blank cells under the new columns, same as the
# existing blank A2/B2 cells.
$ws.Range("A2:B2").Copy($ws.Range("C2"))
$ws.Range("A2:B2").Copy($ws.Range("D2"))
$ws.Range("A2:B2").Copy($ws.Range("E2"))
